$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 09:48:57"
$ws1.Range("A3").Value = "Total filas: 114"

$ws1.Cells.Item(97,2).Value = "09:48:46"
$ws1.Cells.Item(97,3).Value = "09:51"
$ws1.Cells.Item(97,4).Value = "10_OLMOS"
$ws1.Cells.Item(97,5).Value = 3
$ws1.Cells.Item(97,6).Value = "LP1912"
$ws1.Cells.Item(97,7).Value = "30/12/2025"
$ws1.Cells.Item(98,2).Value = "09:48:46"
$ws1.Cells.Item(98,3).Value = "09:52"
$ws1.Cells.Item(98,4).Value = "15_ABASTO"
$ws1.Cells.Item(98,5).Value = 4
$ws1.Cells.Item(98,6).Value = "LP1912"
$ws1.Cells.Item(98,7).Value = "30/12/2025"
$ws1.Cells.Item(99,2).Value = "09:48:46"
$ws1.Cells.Item(99,3).Value = "10:03"
$ws1.Cells.Item(99,4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(99,5).Value = 15
$ws1.Cells.Item(99,6).Value = "LP1912"
$ws1.Cells.Item(99,7).Value = "30/12/2025"
$ws1.Cells.Item(100,2).Value = "09:48:46"
$ws1.Cells.Item(100,3).Value = "10:04"
$ws1.Cells.Item(100,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(100,5).Value = 16
$ws1.Cells.Item(100,6).Value = "LP1912"
$ws1.Cells.Item(100,7).Value = "30/12/2025"
$ws1.Cells.Item(101,2).Value = "09:48:46"
$ws1.Cells.Item(101,3).Value = "10:10"
$ws1.Cells.Item(101,4).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(101,5).Value = 22
$ws1.Cells.Item(101,6).Value = "LP1912"
$ws1.Cells.Item(101,7).Value = "30/12/2025"
$ws1.Cells.Item(102,2).Value = "09:48:46"
$ws1.Cells.Item(102,3).Value = "10:12"
$ws1.Cells.Item(102,4).Value = "15_ABASTO"
$ws1.Cells.Item(102,5).Value = 24
$ws1.Cells.Item(102,6).Value = "LP1912"
$ws1.Cells.Item(102,7).Value = "30/12/2025"
$ws1.Cells.Item(103,2).Value = "09:48:46"
$ws1.Cells.Item(103,3).Value = "10:13"
$ws1.Cells.Item(103,4).Value = "10_OLMOS"
$ws1.Cells.Item(103,5).Value = 25
$ws1.Cells.Item(103,6).Value = "LP1912"
$ws1.Cells.Item(103,7).Value = "30/12/2025"
$ws1.Cells.Item(104,2).Value = "09:48:46"
$ws1.Cells.Item(104,3).Value = "10:21"
$ws1.Cells.Item(104,4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(104,5).Value = 33
$ws1.Cells.Item(104,6).Value = "LP1912"
$ws1.Cells.Item(104,7).Value = "30/12/2025"
$ws1.Cells.Item(105,2).Value = "09:48:46"
$ws1.Cells.Item(105,3).Value = "10:22"
$ws1.Cells.Item(105,4).Value = "17_ROMERO"
$ws1.Cells.Item(105,5).Value = 34
$ws1.Cells.Item(105,6).Value = "LP1912"
$ws1.Cells.Item(105,7).Value = "30/12/2025"
$ws1.Cells.Item(106,2).Value = "09:48:46"
$ws1.Cells.Item(106,3).Value = "10:23"
$ws1.Cells.Item(106,4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(106,5).Value = 35
$ws1.Cells.Item(106,6).Value = "LP1912"
$ws1.Cells.Item(106,7).Value = "30/12/2025"
$ws1.Cells.Item(107,2).Value = "09:48:46"
$ws1.Cells.Item(107,3).Value = "10:26"
$ws1.Cells.Item(107,4).Value = "215A_EL PATO"
$ws1.Cells.Item(107,5).Value = 38
$ws1.Cells.Item(107,6).Value = "LP1912"
$ws1.Cells.Item(107,7).Value = "30/12/2025"
$ws1.Cells.Item(108,2).Value = "09:48:46"
$ws1.Cells.Item(108,3).Value = "10:34"
$ws1.Cells.Item(108,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(108,5).Value = 46
$ws1.Cells.Item(108,6).Value = "LP1912"
$ws1.Cells.Item(108,7).Value = "30/12/2025"
$ws1.Cells.Item(109,2).Value = "09:48:46"
$ws1.Cells.Item(109,3).Value = "10:41"
$ws1.Cells.Item(109,4).Value = "17_ROMERO"
$ws1.Cells.Item(109,5).Value = 53
$ws1.Cells.Item(109,6).Value = "LP1912"
$ws1.Cells.Item(109,7).Value = "30/12/2025"
$ws1.Cells.Item(110,2).Value = "09:48:46"
$ws1.Cells.Item(110,3).Value = "10:43"
$ws1.Cells.Item(110,4).Value = "14_ABASTO"
$ws1.Cells.Item(110,5).Value = 55
$ws1.Cells.Item(110,6).Value = "LP1912"
$ws1.Cells.Item(110,7).Value = "30/12/2025"
$ws1.Cells.Item(111,2).Value = "09:48:46"
$ws1.Cells.Item(111,3).Value = "10:56"
$ws1.Cells.Item(111,4).Value = "27_EL RETIRO"
$ws1.Cells.Item(111,5).Value = 68
$ws1.Cells.Item(111,6).Value = "LP1912"
$ws1.Cells.Item(111,7).Value = "30/12/2025"
$ws1.Cells.Item(112,2).Value = "09:48:46"
$ws1.Cells.Item(112,3).Value = "11:01"
$ws1.Cells.Item(112,4).Value = "215C_EL PATO"
$ws1.Cells.Item(112,5).Value = 73
$ws1.Cells.Item(112,6).Value = "LP1912"
$ws1.Cells.Item(112,7).Value = "30/12/2025"
$ws1.Cells.Item(113,2).Value = "09:48:46"
$ws1.Cells.Item(113,3).Value = "11:13"
$ws1.Cells.Item(113,4).Value = "10_OLMOS"
$ws1.Cells.Item(113,5).Value = 85
$ws1.Cells.Item(113,6).Value = "LP1912"
$ws1.Cells.Item(113,7).Value = "30/12/2025"
$ws1.Cells.Item(114,2).Value = "09:48:46"
$ws1.Cells.Item(114,3).Value = "11:21"
$ws1.Cells.Item(114,4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(114,5).Value = 93
$ws1.Cells.Item(114,6).Value = "LP1912"
$ws1.Cells.Item(114,7).Value = "30/12/2025"
$ws1.Cells.Item(115,2).Value = "09:48:46"
$ws1.Cells.Item(115,3).Value = "11:22"
$ws1.Cells.Item(115,4).Value = "15_ABASTO"
$ws1.Cells.Item(115,5).Value = 94
$ws1.Cells.Item(115,6).Value = "LP1912"
$ws1.Cells.Item(115,7).Value = "30/12/2025"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 09:48:57"
$ws2.Range("A3").Value = "Total filas: 15"

$ws2.Cells.Item(15,2).Value = "30/12/2025"
$ws2.Cells.Item(15,3).Value = "09:48:46"
$ws2.Cells.Item(15,4).Value = "10:26"
$ws2.Cells.Item(15,5).Value = "215A_EL PATO"
$ws2.Cells.Item(15,6).Value = 38
$ws2.Cells.Item(15,7).Value = "LP1912"
$ws2.Cells.Item(16,2).Value = "30/12/2025"
$ws2.Cells.Item(16,3).Value = "09:48:46"
$ws2.Cells.Item(16,4).Value = "11:01"
$ws2.Cells.Item(16,5).Value = "215C_EL PATO"
$ws2.Cells.Item(16,6).Value = 73
$ws2.Cells.Item(16,7).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 30/12/2025 09:48:57"
$ws3.Range("A3").Value = "Total filas: 16"

$ws3.Cells.Item(15,2).Value = "30/12/2025"
$ws3.Cells.Item(15,3).Value = "09:48:57"
$ws3.Cells.Item(15,4).Value = "10:02"
$ws3.Cells.Item(15,5).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(15,6).Value = 14
$ws3.Cells.Item(15,7).Value = "L6173"
$ws3.Cells.Item(16,2).Value = "30/12/2025"
$ws3.Cells.Item(16,3).Value = "09:48:57"
$ws3.Cells.Item(16,4).Value = "10:54"
$ws3.Cells.Item(16,5).Value = "215A_LA PLATA"
$ws3.Cells.Item(16,6).Value = 66
$ws3.Cells.Item(16,7).Value = "L6173"
$ws3.Cells.Item(17,2).Value = "30/12/2025"
$ws3.Cells.Item(17,3).Value = "09:48:52"
$ws3.Cells.Item(17,4).Value = "11:13"
$ws3.Cells.Item(17,5).Value = "215C_LA PLATA"
$ws3.Cells.Item(17,6).Value = 85
$ws3.Cells.Item(17,7).Value = "L6203"
